# Add 2022-Q3 data.
#
# Before: sheet "总计" (totals) + sheet "2022-Q2" (fund holdings for 2022-Q2).
# After:  sheet "总计" (totals, +1 row) + sheet "2022-Q3" (new fund holdings)
#         + sheet "2022-Q2" (the original fund-holdings sheet, re-created
#         under a fresh tab so the old Q2 numbers are preserved).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. The existing "2022-Q2" worksheet becomes the new "2022-Q3" worksheet
#    (it keeps its sheetId / look - only name + data change).
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q2")
$q3.Name = "2022-Q3"

$q3.Range("B2").Value = "'000259"
$q3.Range("C2").Value = "农银区间收益混合"
$q3.Range("D2").Value = "'4.20"
$q3.Range("E2").Value = "'69.93"
$q3.Range("F2").Value = "'2.61"
$q3.Range("G2").Value = "'0.1096"
$q3.Range("H2").Value = 3

$q3.Range("B3").Value = "'005638"
$q3.Range("C3").Value = "农银汇理量化智慧动力混合"
$q3.Range("D3").Value = "'0.60"
$q3.Range("E3").Value = "'88.21"
$q3.Range("F3").Value = "'4.41"
$q3.Range("G3").Value = "'0.0265"
$q3.Range("H3").Value = 1

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "'562530"
$q3.Range("C4").Value = "华夏中证智选1000价值稳健策略ETF"
$q3.Range("D4").Value = "'0.54"
$q3.Range("E4").Value = "'94.32"
$q3.Range("F4").Value = "'0.89"
$q3.Range("G4").Value = "'0.0048"
$q3.Range("H4").Value = 9

# Carry the index-column formatting (bold + border, same as rows 2/3) down
# onto the new row 4.
$q3.Range("A2").Copy()
$q3.Range("A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Re-create a standalone "2022-Q2" sheet (right after "2022-Q3") holding
#    the figures the original "2022-Q2" sheet used to contain.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Add($null, $q3)
$q2.Name = "2022-Q2"

$q2.Range("B1").Value = "基金代码"
$q2.Range("C1").Value = "基金名称"
$q2.Range("D1").Value = "基金规模"
$q2.Range("E1").Value = "股票总仓位"
$q2.Range("F1").Value = "仓位占比"
$q2.Range("G1").Value = "持有市值(亿元)"
$q2.Range("H1").Value = "仓位排名"

$q2.Range("A2").Value = 0
$q2.Range("B2").Value = "'000259"
$q2.Range("C2").Value = "农银区间收益混合"
$q2.Range("D2").Value = "'4.43"
$q2.Range("E2").Value = "'67.87"
$q2.Range("F2").Value = "'1.72"
$q2.Range("G2").Value = "'0.0762"
$q2.Range("H2").Value = 7

$q2.Range("A3").Value = 1
$q2.Range("B3").Value = "'005638"
$q2.Range("C3").Value = "农银汇理量化智慧动力混合"
$q2.Range("D3").Value = "'0.71"
$q2.Range("E3").Value = "'86.53"
$q2.Range("F3").Value = "'2.21"
$q2.Range("G3").Value = "'0.0157"
$q2.Range("H3").Value = 6

# Copy over the header / index-column look (bold + border) from "总计".
$wb.Worksheets.Item("总计").Range("B1").Copy()
$q2.Range("B1:H1").PasteSpecial(-4122)
$wb.Worksheets.Item("总计").Range("A2").Copy()
$q2.Range("A2:A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Update the "总计" (totals) sheet: insert a 2022-Q3 row above the
#    existing 2022-Q2 row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Range("A3").EntireRow.Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.14

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.09

Write-Output "done"
